$wb = $excel.ActiveWorkbook

# ----- Overview sheet: row 3 (b.md) now reports "Ready for handoff" -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-23 22:35:12"

# ----- zh-cn sheet: row 3 (b.md) handoff details -----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps "False" stored as text (matches the other
# Content-Duplicate cells) instead of Excel's automatic Boolean coercion;
# resetting the Style afterwards drops the temporary quote-prefix format.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-23 22:34:59"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2c928257bb94f5696e54fbfd550d53aef159569/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.2

# ----- de-de sheet: row 3 (b.md) handoff details -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-23 22:35:12"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b7ee68aa88eae23c4bc0ec3e6aec3a82a48898eb/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f2c928257bb94f5696e54fbfd550d53aef159569/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.2
